$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "ValidLogin"

# Set data row first (so shared strings pick up Admin/admin123 before Username/Password)
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"

# Set header row
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"

# Select entire row 7 (A7:XFD7), active cell A7
$ws.Rows.Item(7).Select()
